$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert new record at row 434 ---
# Every existing record from row 434 downward shifts down by one row.
$ws.Rows.Item(434).Insert()

$ws.Cells.Item(434,1).Value  = 10
$ws.Cells.Item(434,2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(434,3).Value  = "La Araucanía"
$ws.Cells.Item(434,4).Value  = 44748
$ws.Cells.Item(434,5).Value  = 9
$ws.Cells.Item(434,6).Value  = "Fruta"
$ws.Cells.Item(434,7).Value  = 100108
$ws.Cells.Item(434,8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(434,9).Value  = 100108005
$ws.Cells.Item(434,10).Value = "Piña"
$ws.Cells.Item(434,11).Value = "Caramelo"
$ws.Cells.Item(434,12).Value = "Segunda"
$ws.Cells.Item(434,13).Value = 90
$ws.Cells.Item(434,14).Value = 23000
$ws.Cells.Item(434,15).Value = 23000
$ws.Cells.Item(434,16).Value = 23000
$ws.Cells.Item(434,17).Value = "`$/caja 14 unidades"
$ws.Cells.Item(434,18).Value = "Ecuador"
$ws.Cells.Item(434,19).Value = 1643
$ws.Cells.Item(434,20).Value = 14

# --- Insert a second new record at row 491 ---
# (After the first insert, former row 491 now sits at row 492; inserting here
# pushes it, and the row after it, down one more to rows 492 and 493.)
$ws.Rows.Item(491).Insert()

$ws.Cells.Item(491,1).Value  = 10
$ws.Cells.Item(491,2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(491,3).Value  = "La Araucanía"
$ws.Cells.Item(491,4).Value  = 44519
$ws.Cells.Item(491,5).Value  = 9
$ws.Cells.Item(491,6).Value  = "Fruta"
$ws.Cells.Item(491,7).Value  = 100108
$ws.Cells.Item(491,8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(491,9).Value  = 100108005
$ws.Cells.Item(491,10).Value = "Piña"
$ws.Cells.Item(491,11).Value = "Caramelo"
$ws.Cells.Item(491,12).Value = "Segunda"
$ws.Cells.Item(491,13).Value = 60
$ws.Cells.Item(491,14).Value = 20000
$ws.Cells.Item(491,15).Value = 20000
$ws.Cells.Item(491,16).Value = 20000
$ws.Cells.Item(491,17).Value = "`$/caja 14 unidades"
$ws.Cells.Item(491,18).Value = "Ecuador"
$ws.Cells.Item(491,19).Value = 1429
$ws.Cells.Item(491,20).Value = 14
